$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 171227
$ws.Range("C4").Value = 162034
$ws.Range("C5").Value = 9194
$ws.Range("C8").Value = 65.95
